$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 05:42"

# Update country rows whose label and/or daily figures changed
$ws.Cells.Item(22, 1).Value = "Brasil"
$ws.Cells.Item(22, 2).Value = 2985
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 6
$ws.Cells.Item(22, 5).Value = 2902
$ws.Cells.Item(22, 6).Value = 296
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 77

$ws.Cells.Item(25, 1).Value = "Malasia"
$ws.Cells.Item(25, 2).Value = 2031
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 215
$ws.Cells.Item(25, 5).Value = 1792
$ws.Cells.Item(25, 6).Value = 45
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 24

$ws.Cells.Item(44, 1).Value = "India"
$ws.Cells.Item(44, 2).Value = 733
$ws.Cells.Item(44, 3).Value = 6
$ws.Cells.Item(44, 4).Value = 66
$ws.Cells.Item(44, 5).Value = 647
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 20

$ws.Cells.Item(95, 1).Value = "Kazajistan"
$ws.Cells.Item(95, 2).Value = 121
$ws.Cells.Item(95, 3).Value = 8
$ws.Cells.Item(95, 4).Value = 2
$ws.Cells.Item(95, 5).Value = 118
$ws.Cells.Item(95, 6).Value = 1
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 1

$ws.Cells.Item(96, 1).Value = "Brunei"
$ws.Cells.Item(96, 2).Value = 114
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 5
$ws.Cells.Item(96, 5).Value = 109
$ws.Cells.Item(96, 6).Value = 1
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 0

$ws.Cells.Item(101, 1).Value = "Camboya"
$ws.Cells.Item(101, 2).Value = 98
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 11
$ws.Cells.Item(101, 5).Value = 87
$ws.Cells.Item(101, 6).Value = 1
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 0

$ws.Cells.Item(119, 1).Value = "Paraguay"
$ws.Cells.Item(119, 2).Value = 52
$ws.Cells.Item(119, 3).Value = 11
$ws.Cells.Item(119, 4).Value = 1
$ws.Cells.Item(119, 5).Value = 48
$ws.Cells.Item(119, 6).Value = 1
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 3

$ws.Cells.Item(120, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(120, 2).Value = 51
$ws.Cells.Item(120, 3).Value = 0
$ws.Cells.Item(120, 4).Value = 0
$ws.Cells.Item(120, 5).Value = 48
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 3

$ws.Cells.Item(121, 1).Value = "Ruanda"
$ws.Cells.Item(121, 2).Value = 50
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(121, 4).Value = 0
$ws.Cells.Item(121, 5).Value = 50
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 0

$ws.Cells.Item(122, 1).Value = "Kirguistan"
$ws.Cells.Item(122, 2).Value = 44
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 44
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 0

$ws.Cells.Item(123, 1).Value = "Banglades"
$ws.Cells.Item(123, 2).Value = 44
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 4).Value = 11
$ws.Cells.Item(123, 5).Value = 28
$ws.Cells.Item(123, 6).Value = 1
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 5

$ws.Cells.Item(171, 1).Value = "Montserrat"
$ws.Cells.Item(171, 2).Value = 5
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 5).Value = 5
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 0

$ws.Cells.Item(173, 1).Value = "Fiyi"
$ws.Cells.Item(173, 2).Value = 5
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 5
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

$ws.Cells.Item(175, 1).Value = "Cabo Verde"
$ws.Cells.Item(175, 2).Value = 5
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 5).Value = 4
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 1

$ws.Cells.Item(176, 1).Value = "Guyana"
$ws.Cells.Item(176, 2).Value = 5
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 5).Value = 4
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 1

$ws.Cells.Item(182, 1).Value = "Mauritania"
$ws.Cells.Item(182, 2).Value = 3
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 3
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0

$ws.Cells.Item(183, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(183, 2).Value = 3
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 3
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

$ws.Cells.Item(184, 1).Value = "Liberia"
$ws.Cells.Item(184, 2).Value = 3
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 3
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

$ws.Cells.Item(185, 1).Value = "Republica del Chad"
$ws.Cells.Item(185, 2).Value = 3
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 3
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

$ws.Cells.Item(186, 1).Value = "San Bartolome"
$ws.Cells.Item(186, 2).Value = 3
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 3
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

$ws.Cells.Item(187, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(187, 2).Value = 3
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 3
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

$ws.Cells.Item(188, 1).Value = "Butan"
$ws.Cells.Item(188, 2).Value = 3
$ws.Cells.Item(188, 3).Value = 1
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 5).Value = 3
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

$ws.Cells.Item(189, 1).Value = "Santa Lucia"
$ws.Cells.Item(189, 2).Value = 3
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 1
$ws.Cells.Item(189, 5).Value = 2
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 0

$ws.Cells.Item(190, 1).Value = "Nepal"
$ws.Cells.Item(190, 2).Value = 3
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 1
$ws.Cells.Item(190, 5).Value = 2
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

$ws.Cells.Item(191, 1).Value = "Gambia"
$ws.Cells.Item(191, 2).Value = 3
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 2
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 1

$ws.Cells.Item(193, 1).Value = "Zimbabue"
$ws.Cells.Item(193, 2).Value = 3
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 0
$ws.Cells.Item(193, 5).Value = 2
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 1

$ws.Cells.Item(194, 1).Value = "Anguila"
$ws.Cells.Item(194, 2).Value = 2
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 0
$ws.Cells.Item(194, 5).Value = 2
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0

